$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 0.3300724615171935
$ws.Range("B4").Value = 0.3939180674855193
$ws.Range("C4").Value = 0.28315744431636
$ws.Range("D4").Value = 0.3572073749129407
$ws.Range("E4").Value = 0.3017672444927606
$ws.Range("F4").Value = 7.009940117459498
$ws.Range("G4").Value = 9.14049226783138
$ws.Range("H4").Value = 5.44428906282274
$ws.Range("I4").Value = 7.916869781435245
$ws.Range("J4").Value = 6.061858984930184
$ws.Range("K4").Value = 2.313796406217068
$ws.Range("L4").Value = 3.157439112056929
$ws.Range("M4").Value = 1.704766793562185
$ws.Range("N4").Value = 2.666677956112652
$ws.Range("O4").Value = 1.945404950478435
$ws.Range("P4").Value = 65.082931
$ws.Range("Q4").Value = 90.80950833653185
$ws.Range("R4").Value = 47.4482028160903
$ws.Range("S4").Value = 74.37772062264614
$ws.Range("T4").Value = 54.8171366293158
$ws.Range("U4").Value = 0.2015597452512233
$ws.Range("V4").Value = 0.3171580079148343
$ws.Range("W4").Value = 0.09954007807863229
$ws.Range("X4").Value = 0.2604303568272447
$ws.Range("Y4").Value = 0.1428659843637737
$ws.Range("Z4").Value = 0.8481391630065199
$ws.Range("AA4").Value = 0.9544990526690307
$ws.Range("AB4").Value = 0.6981896653970183
$ws.Range("AC4").Value = 0.9099013217487522
$ws.Range("AD4").Value = 0.7730145981012134
